$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.622.28"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.927.63"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").Value = "'326.52"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "'0.4060"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.08199"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").Value = "'23.74"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.928.58"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "'6.069"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'7.284"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "'91.48"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'17.62"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "29.611.27"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'5.656"
$ws.Range("D23").Value = "'11.94"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "'2.203"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "2.115.82"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "'156.34"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "'20.03"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "'120.77"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("D32").Value = "'0.09609"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'5.627"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "'3.559"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "'0.06543"
$ws.Range("E36").Value = "  +6.80%  "
$ws.Range("D37").Value = "'0.02278"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'1.213"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").Value = "'0.5929"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "'1.011"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.852"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1843"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.484"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.244"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.38"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.07546"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5548"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.961"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'118.23"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").Value = "'2.433"
$ws.Range("E51").Value = "  +0.37%  "
